$wb = $excel.ActiveWorkbook

# --- Sheet "Statistics": simulation speed/density samples ---
$ws1 = $wb.Worksheets.Item("Statistics")

$statsData = @(
    ,@("2024-07-30 22:45:38", 100.0804194481402, 8)
    ,@("2024-07-30 22:45:40", 97.73670682793352, 16)
    ,@("2024-07-30 22:45:42", 97.31738825339654, 24)
    ,@("2024-07-30 22:45:44", 95.23234230223129, 30)
    ,@("2024-07-30 22:45:46", 92.6567857709139, 34)
    ,@("2024-07-30 22:45:48", 89.3488440283692, 37)
    ,@("2024-07-30 22:45:50", 89.14872053614005, 38)
    ,@("2024-07-30 22:45:52", 86.01621069767486, 36)
    ,@("2024-07-30 22:45:54", 82.38666017015372, 38)
    ,@("2024-07-30 22:45:56", 79.62855465991544, 38)
    ,@("2024-07-30 22:45:58", 78.35310442587, 40)
    ,@("2024-07-30 22:46:00", 79.10140228330833, 39)
    ,@("2024-07-30 22:46:02", 76.36352199368717, 45)
    ,@("2024-07-30 22:46:04", 74.95296951464908, 45)
    ,@("2024-07-30 22:46:06", 78.40997691535971, 44)
    ,@("2024-07-30 22:46:08", 80.40935495769017, 40)
)

$r = 2
foreach ($row in $statsData) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# --- Sheet "Accidents": lane-switching crash log ---
$ws2 = $wb.Worksheets.Item("Accidents")

$accidentData = @(
    ,@("2024-07-30 22:46:03", "Car and Car", "30.33 and 55.83")
    ,@("2024-07-30 22:46:04", "Car and Car", "64.24 and 0.00")
)

$r = 2
foreach ($row in $accidentData) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
